# Update the cryptos list with the latest scraped price/volume figures.
# Numeric-looking price strings are prefixed with a literal leading
# apostrophe so Excel keeps them as text (matching the original
# inline-string cell type) instead of silently re-parsing them as
# numbers and dropping significant trailing/leading zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.400.51'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '1.688.49'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').Value = '  +0.87%  '
$ws.Range('D5').Value = '''218.62'
$ws.Range('D6').Value = '''0.5463'
$ws.Range('E6').Value = '  +6.96%  '
$ws.Range('E7').Value = '  +0.81%  '
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('D9').Value = '''0.06460'
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').Value = '''0.07690'
$ws.Range('E11').Value = '  +3.33%  '
$ws.Range('D12').Value = '1.690.09'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').Value = '''4.533'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').Value = '''0.5808'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '''0.000008344'
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('D16').Value = '''65.17'
$ws.Range('E16').Value = '  +1.67%  '
$ws.Range('D17').Value = '26.451.02'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').Value = '''4.957'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').Value = '''189.95'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '''6.220'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '''1.013'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').Value = '''150.41'
$ws.Range('E24').Value = '  +3.93%  '
$ws.Range('D25').Value = '''0.1303'
$ws.Range('E25').Value = '  +6.24%  '
$ws.Range('E26').Value = '  +3.74%  '
$ws.Range('D27').Value = '''15.70'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').Value = '''0.06351'
$ws.Range('E28').Value = '  -4.44%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''1.413'
$ws.Range('E29').Value = '  +5.46%  '
$ws.Range('D30').Value = '''1.328'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').Value = '''3.580'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').Value = '''3.570'
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('D33').Value = '''1.675'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('E34').Value = '  +2.38%  '
$ws.Range('D35').Value = '''0.6208'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '''2.417'
$ws.Range('E36').Value = '  +2.09%  '
$ws.Range('D37').Value = '''2.723'
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('D38').Value = '''6.211'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').Value = '1.117.17'
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('D40').Value = '''0.01636'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').Value = '''0.8807'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').Value = '''101.09'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '1.842.68'
$ws.Range('E44').Value = '  +1.49%  '
$ws.Range('E45').Value = '  -5.16%  '
$ws.Range('D46').Value = '''57.27'
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').Value = '''8.221'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').Value = '''1.009'
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').Value = '''0.05277'
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('D50').Value = '''0.4304'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '''6.055'
$ws.Range('E51').Value = '  +1.22%  '
